# The sheet grows from a single "read a double value from an xlsx" demo
# column (A1:A5, dimension A1:A5) into a 3-column dataset of CO2 values
# for three car brands (mercedes / audi / bmw), still with the numbers
# kept as text (matching the original file's own convention of storing
# "5.6" etc. as shared-string text rather than numeric cells) so that a
# chosen column index can be read back as a string, per the commit
# message ("csv en excel kan gelezen worden met gekozen kolom index").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "co2 mercedes"
$ws.Range("B1").Value = "co2 audi"
$ws.Range("C1").Value = "co2 bmw"

# Data rows — written as formulas yielding text so Excel doesn't coerce
# the numeric-looking strings ("5.6", "7.9", ...) into numbers.
$ws.Range("A2").Formula = "=""5.6"""
$ws.Range("B2").Formula = "=""7.9"""
$ws.Range("C2").Formula = "=""9.5"""

$ws.Range("A3").Formula = "=""3.8"""
$ws.Range("B3").Formula = "=""8.4"""
$ws.Range("C3").Formula = "=""5.7"""

$ws.Range("A4").Formula = "=""9.4"""
$ws.Range("B4").Formula = "=""1.6"""
$ws.Range("C4").Formula = "=""10.4"""

$ws.Range("A5").Formula = "=""15.6"""
$ws.Range("B5").Formula = "=""4.8"""
$ws.Range("C5").Formula = "=""16.4"""

# Convert the formulas to plain text values in place (copy + paste
# special values only) so the cells end up as ordinary shared-string
# text cells, not formulas, and without touching any cell's style.
$full = $ws.Range("A1:C5")
$full.Copy()
$full.PasteSpecial(-4163)

# Match the author's final selection.
$ws.Range("F5").Select()
